$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template used to generate a select-list column for each of
# "is_locked" and "is_enabled" dict fields (columns D and E). Those two
# columns are being dropped; the remaining "order_by" / "rem" columns
# shift left into D/E, and the trailing F/G columns are removed.
$ws.Range("D1").Value = "<%=comment.order_by%>"
$ws.Range("E1").Value = "<%=comment.rem%>"
$ws.Range("F1:G1").ClearContents() | Out-Null
